# Correct the CV keyword list in column E (CV KEYWORDS) on the active sheet.
# The keyword extraction logic had a bug when matching keywords to the
# correct experience/section entry; this resulted in wrong counts for some
# keywords and a different ordering of the aggregated keyword list.
# This script rewrites the affected cells with the corrected keyword list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "python : 8"
$ws.Range("E20").Value = "reports : 1"
$ws.Range("E21").Value = "qlik : 4"
$ws.Range("E22").Value = "conversion : 1"
$ws.Range("E23").Value = "sql : 1"
$ws.Range("E24").Value = "ssis : 1"
$ws.Range("E25").Value = "etl : 1"
$ws.Range("E26").Value = "selenium : 2"
$ws.Range("E27").Value = "sql queries : 1"
$ws.Range("E28").Value = "queries : 1"
$ws.Range("E29").Value = "stored procedures : 1"
$ws.Range("E30").Value = "bi : 1"
$ws.Range("E31").Value = "developer : 3"
$ws.Range("E38").Value = "data mining : 3"
$ws.Range("E40").Value = "tensorflow : 3"
$ws.Range("E47").Value = "hybrid : 1"
$ws.Range("E48").Value = "design : 1"
$ws.Range("E50").Value = "tools : 2"
$ws.Range("E51").Value = "docker : 2"
$ws.Range("E52").Value = "ansible : 2"
$ws.Range("E53").Value = "cloudformation : 1"
$ws.Range("E54").Value = "azure : 2"
$ws.Range("E55").Value = "amazon : 1"
$ws.Range("E56").Value = "web services : 1"
$ws.Range("E57").Value = "microsoft azure : 1"
$ws.Range("E58").Value = "databases : 1"
$ws.Range("E59").Value = "redshift : 2"
$ws.Range("E60").Value = "mysql : 2"
$ws.Range("E61").Value = "unix shell : 2"
$ws.Range("E62").Value = "shell : 1"
$ws.Range("E63").Value = "java : 1"
